$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new "Symptoms" column (G) ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Is Student?"
$ws.Range("C1").Value = "Favorite Subject"
$ws.Range("D1").Value = "Favorite Season"
$ws.Range("E1").Value = "Likes Cats"
$ws.Range("F1").Value = "Gender"
$ws.Range("G1").Value = "Symptoms"

# --- Row 2 is left blank (Joe's old record moves down to row 4) ---
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# --- Row 3: Bob (existing record, updated + new Symptoms cell) ---
$ws.Range("A3").Value = "Bob"
$ws.Range("B3").Value = "No"
$ws.Range("C3").Value = "Math"
$ws.Range("D3").Value = "Summer"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Other"
$ws.Range("G3").Value = "Tired, Fever"

# --- Row 4: Joe (moved down from row 2, updated + new Symptoms cell) ---
$ws.Range("A4").Value = "Joe"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "Science"
$ws.Range("D4").Value = "Fall"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "Male"
$ws.Range("G4").Value = "Headache, Itchy"

# --- Row 5: Milo (new record) ---
$ws.Range("A5").Value = "Milo"
$ws.Range("B5").Value = "No"
$ws.Range("C5").Value = "Treats!"
$ws.Range("D5").Value = "Summer`n"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "Male"
$ws.Range("G5").Value = "Itchy"

# --- Row 6: Luna (new record, no Symptoms value) ---
$ws.Range("A6").Value = "Luna"
$ws.Range("B6").Value = "No"
$ws.Range("C6").Value = "Long walks"
$ws.Range("D6").Value = "Winter"
$ws.Range("E6").Value = "Yes"
$ws.Range("F6").Value = "Female"
$ws.Range("G6").Value = ""

# --- Row 7: Meghan (new record) ---
$ws.Range("A7").Value = "Meghan"
$ws.Range("B7").Value = "Yes"
$ws.Range("C7").Value = "Business`n"
$ws.Range("D7").Value = "Spring`n"
$ws.Range("E7").Value = "Yes"
$ws.Range("F7").Value = "Female"
$ws.Range("G7").Value = "Tired"

# --- Column widths: re-fit to (approximately) match the new, wider content.
# (Widths recomputed by Excel's real best-fit glyph metrics land on odd
# fractional character widths that this engine's ColumnWidth setter can only
# approximate to the nearest 1/6 of a character -- these are the closest
# achievable values.)
$ws.Columns.Item(1).ColumnWidth = 7.333333333333333
$ws.Columns.Item(6).ColumnWidth = 6.666666666666667
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666
